# Update crypto price/volume data per latest scrape (GitHub Actions)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = "'63.268.88"
$ws.Range('E2').Value = "'  +0.74%  "
$ws.Range('D3').Value = "'2.551.06"
$ws.Range('E3').Value = "'  +3.70%  "
$ws.Range('E4').Value = "'  -0.05%  "
$ws.Range('D5').Value = "'568.89"
$ws.Range('E5').Value = "'  +1.06%  "
$ws.Range('D6').Value = "'147.63"
$ws.Range('E6').Value = "'  +4.38%  "
$ws.Range('E7').Value = "'  -0.05%  "
$ws.Range('D8').Value = "'0.588"
$ws.Range('E8').Value = "'  +0.23%  "
$ws.Range('D9').Value = "'2.549.43"
$ws.Range('E9').Value = "'  +3.74%  "
$ws.Range('E10').Value = "'  +0.91%  "
$ws.Range('E11').Value = "'  -1.61%  "
$ws.Range('E12').Value = "'  +0.66%  "
$ws.Range('D13').Value = "'0.353"
$ws.Range('E13').Value = "'  +0.55%  "
$ws.Range('D14').Value = "'27.58"
$ws.Range('E14').Value = "'  +5.38%  "
$ws.Range('D15').Value = "'3.004.69"
$ws.Range('E15').Value = "'  +3.59%  "
$ws.Range('D16').Value = "'63.178.50"
$ws.Range('E16').Value = "'  +0.69%  "
$ws.Range('E17').Value = "'  +2.47%  "
$ws.Range('D18').Value = "'2.547.56"
$ws.Range('E18').Value = "'  +3.42%  "
$ws.Range('E19').Value = "'  +2.33%  "
$ws.Range('D20').Value = "'335.89"
$ws.Range('E20').Value = "'  -0.90%  "
$ws.Range('D21').Value = "'4.31"
$ws.Range('E21').Value = "'  +1.48%  "
$ws.Range('D22').Value = "'6.78"
$ws.Range('E22').Value = "'  +0.10%  "
$ws.Range('E23').Value = "'  +0.24%  "
$ws.Range('D24').Value = "'65.33"
$ws.Range('E24').Value = "'  -0.05%  "
$ws.Range('B25').Value = "'Fetch.AI"
$ws.Range('C25').Value = "'https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet"
$ws.Range('D25').Value = "'1.64"
$ws.Range('E25').Value = "'  +9.94%  "
$ws.Range('B26').Value = "'Kaspa"
$ws.Range('C26').Value = "'https://coinranking.com/coin/V8GxkwWow+kaspa-kas"
$ws.Range('D26').Value = "'0.170"
$ws.Range('E26').Value = "'  -1.91%  "
$ws.Range('D27').Value = "'1.51"
$ws.Range('E27').Value = "'  +11.10%  "
$ws.Range('D28').Value = "'8.52"
$ws.Range('E28').Value = "'  +6.11%  "
$ws.Range('D29').Value = "'0.999"
$ws.Range('E29').Value = "'  -0.24%  "
$ws.Range('D30').Value = "'7.37"
$ws.Range('E30').Value = "'  +8.40%  "
$ws.Range('D31').Value = "'0.0₃0825"
$ws.Range('E31').Value = "'  +3.71%  "
$ws.Range('D32').Value = "'1.85"
$ws.Range('E32').Value = "'  +0.87%  "
$ws.Range('D33').Value = "'176.56"
$ws.Range('E33').Value = "'  -0.02%  "
$ws.Range('D34').Value = "'1.57"
$ws.Range('E34').Value = "'  +4.58%  "
$ws.Range('D35').Value = "'413.02"
$ws.Range('E35').Value = "'  +13.82%  "
$ws.Range('E36').Value = "'  +1.58%  "
$ws.Range('E37').Value = "'  +1.22%  "
$ws.Range('D38').Value = "'4.41"
$ws.Range('E38').Value = "'  +1.20%  "
$ws.Range('E39').Value = "'  -0.01%  "
$ws.Range('D40').Value = "'1.77"
$ws.Range('E40').Value = "'  +4.57%  "
$ws.Range('D42').Value = "'39.30"
$ws.Range('E42').Value = "'  -3.00%  "
$ws.Range('D43').Value = "'152.90"
$ws.Range('E43').Value = "'  +2.63%  "
$ws.Range('D44').Value = "'3.79"
$ws.Range('E44').Value = "'  +2.83%  "
$ws.Range('D45').Value = "'21.03"
$ws.Range('E45').Value = "'  +2.77%  "
$ws.Range('D46').Value = "'0.606"
$ws.Range('E46').Value = "'  +1.82%  "
$ws.Range('E47').Value = "'  +1.03%  "
$ws.Range('E48').Value = "'  +6.56%  "
$ws.Range('D49').Value = "'0.0524"
$ws.Range('E49').Value = "'  +1.94%  "
$ws.Range('D50').Value = "'18.37"
$ws.Range('E50').Value = "'  +2.92%  "
$ws.Range('D51').Value = "'1.79"
$ws.Range('E51').Value = "'  +3.35%  "
